$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new daily rows (10/05/2021 - 13/05/2021) after the existing
# last row (251), continuing the date series with zero counts.
# New rows mirror the formatting of the preceding data row (style on
# column A carries the date number format).

$dates = 44326, 44327, 44328, 44329
$startRow = 252

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Copy formatting from the row above (keeps date style "s=2" on col A)
    $ws.Range("A" + ($row - 1) + ":D" + ($row - 1)).Copy($ws.Range("A" + $row + ":D" + $row))

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
